$d = $word.ActiveDocument

# 1) "Underfunktion " -> "Brugermål" + " " (as two separate runs)
$d.Content.Find.Execute("Underfunktion ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Brugermål ", 2)

# 2) Merge the split "(Frequency of Occurrence)" run and remove the stray
#    _GoBack bookmark that separated "Occ" and "urrence)".
$d.Content.Find.Execute(" (Frequency of Occurrence)", $true, $false, $false, $false, $false,
                         $true, 1, $false, " (Frequency of Occurrence)", 2)
